# Apply crypto list price/volume updates (and the WrappedBTC/Polkadot row swap)
# to the active worksheet, matching the commit's refreshed scrape data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper: force the written value to be stored as text (matches the original
# inlineStr cell type) even when the text looks numeric (e.g. "1.00", "595.63"),
# then strip the temporary Text number-format so the cell's style is left untouched.
function Set-TextValue([string]$cellRef, [string]$val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue 'D2' '62.927.48'
Set-TextValue 'E2' '  +3.09%  '
Set-TextValue 'D3' '3.026.23'
Set-TextValue 'E3' '  +2.03%  '
Set-TextValue 'D4' '1.00'
Set-TextValue 'E4' '  -0.09%  '
Set-TextValue 'D5' '595.63'
Set-TextValue 'E5' '  +0.67%  '
Set-TextValue 'D6' '152.95'
Set-TextValue 'E6' '  +7.51%  '
Set-TextValue 'E7' '  -0.27%  '
Set-TextValue 'D8' '3.022.45'
Set-TextValue 'E8' '  +2.07%  '
Set-TextValue 'E9' '  +0.44%  '
Set-TextValue 'D10' '6.94'
Set-TextValue 'E10' '  +15.98%  '
Set-TextValue 'E11' '  +2.39%  '
Set-TextValue 'E12' '  +3.03%  '
Set-TextValue 'D13' '0.0000232'
Set-TextValue 'E13' '  +3.47%  '
Set-TextValue 'D14' '35.65'
Set-TextValue 'E14' '  +5.26%  '
Set-TextValue 'D15' '0.126'
Set-TextValue 'E15' '  -0.22%  '
Set-TextValue 'D16' '3.526.42'
Set-TextValue 'E16' '  +1.72%  '
Set-TextValue 'B17' 'WrappedBTC'
Set-TextValue 'C17' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue 'D17' '62.934.32'
Set-TextValue 'E17' '  +3.00%  '
Set-TextValue 'B18' 'Polkadot'
Set-TextValue 'C18' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D18' '7.06'
Set-TextValue 'E18' '  +3.05%  '
Set-TextValue 'D19' '3.024.86'
Set-TextValue 'E19' '  +1.78%  '
Set-TextValue 'D20' '448.54'
Set-TextValue 'E20' '  +0.60%  '
Set-TextValue 'E21' '  +2.23%  '
Set-TextValue 'E22' '  +2.71%  '
Set-TextValue 'E23' '  +3.45%  '
Set-TextValue 'D24' '82.86'
Set-TextValue 'E24' '  +1.80%  '
Set-TextValue 'D25' '11.38'
Set-TextValue 'E25' '  +8.80%  '
Set-TextValue 'D26' '2.30'
Set-TextValue 'E26' '  +6.50%  '
Set-TextValue 'D27' '12.34'
Set-TextValue 'E27' '  +3.74%  '
Set-TextValue 'E28' '  +0.01%  '
Set-TextValue 'D29' '7.47'
Set-TextValue 'E29' '  +4.93%  '
Set-TextValue 'D30' '2.27'
Set-TextValue 'E30' '  +11.58%  '
Set-TextValue 'E31' '  +1.14%  '
Set-TextValue 'E32' '  -0.16%  '
Set-TextValue 'E33' '  +2.59%  '
Set-TextValue 'D34' '0.110'
Set-TextValue 'E34' '  +1.45%  '
Set-TextValue 'D35' '0.0₃0871'
Set-TextValue 'E35' '  +8.50%  '
Set-TextValue 'E36' '  +3.00%  '
Set-TextValue 'D37' '5.87'
Set-TextValue 'E37' '  +2.45%  '
Set-TextValue 'D38' '3.12'
Set-TextValue 'E38' '  +11.00%  '
Set-TextValue 'E39' '  +8.65%  '
Set-TextValue 'D40' '2.09'
Set-TextValue 'E40' '  +3.64%  '
Set-TextValue 'D41' '50.50'
Set-TextValue 'E41' '  +1.12%  '
Set-TextValue 'E42' '  +1.17%  '
Set-TextValue 'E43' '  +16.96%  '
Set-TextValue 'D44' '0.306'
Set-TextValue 'E44' '  +15.20%  '
Set-TextValue 'D45' '391.48'
Set-TextValue 'E45' '  +2.21%  '
Set-TextValue 'D46' '0.0360'
Set-TextValue 'E46' '  +3.15%  '
Set-TextValue 'D47' '2.708.56'
Set-TextValue 'E47' '  +1.16%  '
Set-TextValue 'D48' '133.88'
Set-TextValue 'E48' '  +2.73%  '
Set-TextValue 'D49' '26.69'
Set-TextValue 'E49' '  +15.51%  '
Set-TextValue 'D51' '2.27'
Set-TextValue 'E51' '  +6.64%  '
